# Append the latest daily profit row (run date 2025-11-02) to Sheet1.
# The Date column stores plain text like "MM/DD/YYYY" (matching every
# other row), so a leading apostrophe is used to stop Excel's COM layer
# from auto-converting the literal into a date serial number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 76
$newRow = $lastRow + 1

$ws.Range("A$newRow").Value = "'11/02/2025"
$ws.Range("B$newRow").Value = 10748.17
